$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.730.20'
$ws.Range('E2').Value = '  -1.92%  '
$ws.Range('D3').Value = '1.754.86'
$ws.Range('E3').Value = '  -2.52%  '
$ws.Range('D4').Value = '''1.003'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''324.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.32%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').Value = '''0.4404'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.03%  '
$ws.Range('D8').Value = '''0.3673'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.84%  '
$ws.Range('D9').Value = '''45.31'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '
$ws.Range('D10').Value = '''0.07486'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('D11').Value = '''1.123'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.09%  '
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').Value = '''21.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.05%  '
$ws.Range('D14').Value = '''6.172'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = '''7.269'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.36%  '
$ws.Range('D16').Value = '1.752.78'
$ws.Range('E16').Value = '  -2.67%  '
$ws.Range('D17').Value = '''0.00001070'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.96%  '
$ws.Range('D18').Value = '''88.27'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +8.39%  '
$ws.Range('D19').Value = '''0.06221'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -7.57%  '
$ws.Range('E20').Value = '  -0.13%  '
$ws.Range('D21').Value = '''17.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.50%  '
$ws.Range('D22').Value = '''6.170'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.65%  '
$ws.Range('D23').Value = '''0.5305'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.21%  '
$ws.Range('D24').Value = '27.732.23'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('D25').Value = '''11.66'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.74%  '
$ws.Range('D26').Value = '''2.329'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.91%  '
$ws.Range('D27').Value = '''20.65'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''154.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '''2.365'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = '1.948.64'
$ws.Range('E30').Value = '  -2.90%  '
$ws.Range('D31').Value = '''128.06'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.50%  '
$ws.Range('D32').Value = '''1.222'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.04%  '
$ws.Range('D33').Value = '''5.746'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.42%  '
$ws.Range('D34').Value = '''0.09179'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.34%  '
$ws.Range('D35').Value = '''3.653'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -9.54%  '
$ws.Range('D36').Value = '''12.68'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +5.20%  '
$ws.Range('D37').Value = '''0.02318'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.20%  '
$ws.Range('D38').Value = '''0.2169'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.54%  '
$ws.Range('D39').Value = '''5.109'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.57%  '
$ws.Range('D40').Value = '''0.6495'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('D41').Value = '''0.06125'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.56%  '
$ws.Range('D42').Value = '''1.198'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.09%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '''7.984'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -5.09%  '
$ws.Range('B44').Value = 'WEMIXTOKEN'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''1.418'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.36%  '
$ws.Range('E45').Value = '  -0.18%  '
$ws.Range('D46').Value = '''13.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.5956'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.43%  '
$ws.Range('B48').Value = 'PancakeSwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D48').Value = '''3.755'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.08%  '
$ws.Range('D49').Value = '''126.20'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.50%  '
$ws.Range('D50').Value = '''1.982'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('D51').Value = '''0.06904'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.45%  '
